# Update "想去人数" (F column) counts across the sheets, matching the
# output regenerated at commit 456a3b4 for gh-pages.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 8315
$ws.Range("F7").Value  = 820
$ws.Range("F8").Value  = 640
$ws.Range("F10").Value = 72
$ws.Range("F12").Value = 879
$ws.Range("F13").Value = 3575
$ws.Range("F14").Value = 248
$ws.Range("F15").Value = 133
$ws.Range("F16").Value = 777
$ws.Range("F22").Value = 799
$ws.Range("F23").Value = 1343
$ws.Range("F24").Value = 394
$ws.Range("F25").Value = 288
$ws.Range("F27").Value = 139
$ws.Range("F28").Value = 320
$ws.Range("F29").Value = 46
$ws.Range("F33").Value = 621
$ws.Range("F34").Value = 37
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 61
$ws.Range("F39").Value = 124

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 211

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 211
$ws.Range("F4").Value  = 8315
$ws.Range("F8").Value  = 820
$ws.Range("F9").Value  = 640
$ws.Range("F11").Value = 72
$ws.Range("F13").Value = 879
$ws.Range("F15").Value = 3575
$ws.Range("F16").Value = 248
$ws.Range("F17").Value = 133
$ws.Range("F19").Value = 777
$ws.Range("F27").Value = 799
$ws.Range("F28").Value = 1343
$ws.Range("F29").Value = 394
$ws.Range("F30").Value = 288
$ws.Range("F32").Value = 139
$ws.Range("F34").Value = 320
$ws.Range("F35").Value = 46
$ws.Range("F39").Value = 621
$ws.Range("F40").Value = 37
$ws.Range("F41").Value = 41
$ws.Range("F42").Value = 61
$ws.Range("F45").Value = 124
